$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Implemented select channel to preview" - fill the preview/channel label column (G)
# for every measurement row (rows 4..81) with its channel id "B1".."B78",
# matching the row's position in the table (row 4 -> B1, row 5 -> B2, ...).
for ($i = 4; $i -le 81; $i++) {
    $b = $i - 3
    $ws.Range("G$i").Value = "B$b"
}

# The last data row (81) loses its bottom border: the table is now visually
# closed off by a new trailing row (82) whose only formatting is a top border,
# so the border shows up between row 81 and row 82 instead of right under 81.
$ws.Range("G81").Borders.Item(9).LineStyle = 0

# Build row 82: copy row 81's cell format (font etc.) into G82 so it keeps the
# same font as the rest of the channel column, then strip it down to just a
# top border.
$ws.Range("G4").Copy()
$ws.Range("G82").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("G82").Borders.LineStyle = 0
$ws.Range("G82").Borders.Item(8).Color = 0
$ws.Range("G82").Borders.Item(8).LineStyle = 1
